# Apply the weekly data refresh: shift existing rows 922-1005 down by two
# positions (losing the two oldest records) and append two new trailing
# rows (1006, 1007) that preserve the records pushed past row 1005, while
# inserting two brand-new rows of data at 922-923.
#
# Columns A, B, C, E, F, G, R are constant across this block, so only the
# D, H, I, J, K, L, M, N, O, P, Q columns need to be written per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(922, '45166', 'Zafiro rojo', 'Primera', '300', '13000', '13000', '13000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '867', '15'),
    @(923, '45166', 'Zafiro verde', 'Primera', '300', '15000', '15000', '15000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1000', '15'),
    @(924, '44397', 'Zafiro rojo', 'Primera', '300', '16000', '16000', '16000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1067', '15'),
    @(925, '44397', 'Zafiro verde', 'Primera', '300', '14000', '14000', '14000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '933', '15'),
    @(926, '44336', 'Zafiro rojo', 'Primera', '200', '27000', '27000', '27000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1800', '15'),
    @(927, '44336', 'Zafiro verde', 'Primera', '300', '15000', '15000', '15000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1000', '15'),
    @(928, '44432', 'Zafiro rojo', 'Primera', '200', '38000', '38000', '38000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2533', '15'),
    @(929, '44432', 'Zafiro verde', 'Primera', '200', '35000', '35000', '35000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2333', '15'),
    @(930, '44901', 'Cuatro cascos verde', 'Primera', '400', '13000', '13000', '13000', '$/caja 15 kilos', 'Región del Maule', '867', '15'),
    @(931, '44217', 'Cuatro cascos verde', 'Primera', '400', '7000', '7000', '7000', '$/caja 15 kilos', 'Región del Maule', '467', '15'),
    @(932, '44767', 'Zafiro rojo', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(933, '44767', 'Zafiro verde', 'Primera', '300', '17000', '17000', '17000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1133', '15'),
    @(934, '44599', 'Cuatro cascos verde', 'Primera', '400', '7000', '7000', '7000', '$/caja 15 kilos', 'Región del Maule', '467', '15'),
    @(935, '44659', 'Cuatro cascos rojo', 'Primera', '200', '28000', '28000', '28000', '$/caja 15 kilos', 'Región del Maule', '1867', '15'),
    @(936, '44659', 'Cuatro cascos verde', 'Primera', '200', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(937, '44459', 'Zafiro rojo', 'Primera', '300', '35000', '35000', '35000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2333', '15'),
    @(938, '44459', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(939, '44880', 'Cuatro cascos verde', 'Primera', '200', '19000', '19000', '19000', '$/caja 18 kilos', 'Región del Maule', '1056', '18'),
    @(940, '45117', 'Morrón rojo', 'Primera', '300', '12000', '12000', '12000', '$/caja 20 kilos', 'Provincia del Elquí', '600', '20'),
    @(941, '45117', 'Morrón rojo', 'Segunda', '200', '10000', '10000', '10000', '$/caja 20 kilos', 'Provincia del Elquí', '500', '20'),
    @(942, '45117', 'Zafiro rojo', 'Primera', '500', '12000', '12000', '12000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '800', '15'),
    @(943, '45117', 'Zafiro verde', 'Primera', '500', '11000', '11000', '11000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '733', '15'),
    @(944, '44637', 'Cuatro cascos rojo', 'Primera', '200', '15000', '15000', '15000', '$/caja 15 kilos', 'Región del Maule', '1000', '15'),
    @(945, '44637', 'Cuatro cascos verde', 'Primera', '300', '8000', '8000', '8000', '$/caja 15 kilos', 'Región del Maule', '533', '15'),
    @(946, '44637', 'Zafiro rojo', 'Primera', '200', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(947, '45083', 'Zafiro rojo', 'Primera', '500', '12000', '12000', '12000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '800', '15'),
    @(948, '45083', 'Zafiro verde', 'Primera', '500', '10000', '10000', '10000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '667', '15'),
    @(949, '44362', 'Zafiro rojo', 'Primera', '300', '15000', '15000', '15000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1000', '15'),
    @(950, '44362', 'Zafiro verde', 'Primera', '300', '11000', '11000', '11000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '733', '15'),
    @(951, '45063', 'Zafiro rojo', 'Primera', '200', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(952, '45063', 'Zafiro verde', 'Primera', '300', '13000', '13000', '13000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '867', '15'),
    @(953, '44792', 'Zafiro rojo', 'Primera', '200', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(954, '44792', 'Zafiro verde', 'Primera', '200', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(955, '44893', 'Cuatro cascos verde', 'Primera', '300', '12000', '12000', '12000', '$/caja 18 kilos', 'Región del Maule', '667', '18'),
    @(956, '44557', 'Cuatro cascos verde', 'Primera', '300', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(957, '44516', 'Cuatro cascos verde', 'Primera', '200', '18000', '18000', '18000', '$/caja 15 kilos', 'Región del Maule', '1200', '15'),
    @(958, '44242', 'Cuatro cascos rojo', 'Primera', '200', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(959, '44242', 'Cuatro cascos verde', 'Primera', '400', '5000', '5000', '5000', '$/caja 15 kilos', 'Región del Maule', '333', '15'),
    @(960, '44848', 'Zafiro rojo', 'Primera', '300', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(961, '44848', 'Zafiro verde', 'Primera', '300', '22000', '22000', '22000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1467', '15'),
    @(962, '44635', 'Cuatro cascos rojo', 'Primera', '200', '15000', '15000', '15000', '$/caja 15 kilos', 'Región del Maule', '1000', '15'),
    @(963, '44635', 'Cuatro cascos verde', 'Primera', '300', '8000', '8000', '8000', '$/caja 15 kilos', 'Región del Maule', '533', '15'),
    @(964, '44635', 'Zafiro rojo', 'Primera', '300', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(965, '44678', 'Zafiro verde', 'Primera', '300', '13000', '13000', '13000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '867', '15'),
    @(966, '44194', 'Cuatro cascos verde', 'Primera', '300', '12000', '12000', '12000', '$/caja 15 kilos', 'Región del Maule', '800', '15'),
    @(967, '44237', 'Cuatro cascos rojo', 'Primera', '200', '13000', '13000', '13000', '$/caja 15 kilos', 'Región del Maule', '867', '15'),
    @(968, '44237', 'Cuatro cascos verde', 'Primera', '300', '6000', '6000', '6000', '$/caja 15 kilos', 'Región del Maule', '400', '15'),
    @(969, '44771', 'Zafiro rojo', 'Primera', '200', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(970, '44771', 'Zafiro verde', 'Primera', '300', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(971, '45051', 'Zafiro rojo', 'Primera', '200', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(972, '44383', 'Zafiro rojo', 'Primera', '400', '16000', '16000', '16000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1067', '15'),
    @(973, '45015', 'Morrón rojo', 'Primera', '200', '12000', '12000', '12000', '$/caja 20 kilos', 'Provincia del Elquí', '600', '20'),
    @(974, '45015', 'Zafiro rojo', 'Primera', '200', '15000', '15000', '15000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1000', '15'),
    @(975, '44530', 'Cuatro cascos verde', 'Primera', '200', '15000', '15000', '15000', '$/caja 15 kilos', 'Región del Maule', '1000', '15'),
    @(976, '44784', 'Zafiro rojo', 'Primera', '300', '23000', '23000', '23000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1533', '15'),
    @(977, '44784', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(978, '44433', 'Zafiro rojo', 'Primera', '200', '38000', '38000', '38000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2533', '15'),
    @(979, '44433', 'Zafiro verde', 'Primera', '200', '35000', '35000', '35000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2333', '15'),
    @(980, '44813', 'Zafiro rojo', 'Primera', '300', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(981, '44813', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(982, '44571', 'Cuatro cascos verde', 'Primera', '300', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(983, '44638', 'Cuatro cascos rojo', 'Primera', '200', '15000', '15000', '15000', '$/caja 15 kilos', 'Región del Maule', '1000', '15'),
    @(984, '44638', 'Cuatro cascos verde', 'Primera', '300', '8000', '8000', '8000', '$/caja 15 kilos', 'Región del Maule', '533', '15'),
    @(985, '44572', 'Cuatro cascos verde', 'Primera', '200', '9000', '9000', '9000', '$/caja 15 kilos', 'Región del Maule', '600', '15'),
    @(986, '44271', 'Cuatro cascos rojo', 'Primera', '200', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(987, '44271', 'Cuatro cascos verde', 'Primera', '200', '6000', '6000', '6000', '$/caja 15 kilos', 'Región del Maule', '400', '15'),
    @(988, '44251', 'Cuatro cascos rojo', 'Primera', '200', '9000', '9000', '9000', '$/caja 15 kilos', 'Región del Maule', '600', '15'),
    @(989, '44251', 'Cuatro cascos verde', 'Primera', '300', '5000', '5000', '5000', '$/caja 15 kilos', 'Región del Maule', '333', '15'),
    @(990, '44286', 'Cuatro cascos rojo', 'Primera', '150', '10000', '10000', '10000', '$/caja 15 kilos', 'Región del Maule', '667', '15'),
    @(991, '44286', 'Cuatro cascos verde', 'Primera', '200', '7000', '7000', '7000', '$/caja 15 kilos', 'Región del Maule', '467', '15'),
    @(992, '44711', 'Zafiro rojo', 'Primera', '300', '45000', '45000', '45000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '3000', '15'),
    @(993, '44711', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(994, '44343', 'Morrón rojo', 'Primera', '200', '17000', '17000', '17000', '$/caja 18 kilos', 'Provincia del Elquí', '944', '18'),
    @(995, '44343', 'Zafiro rojo', 'Primera', '200', '27000', '27000', '27000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1800', '15'),
    @(996, '44343', 'Zafiro verde', 'Primera', '400', '14000', '14000', '14000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '933', '15'),
    @(997, '44754', 'Zafiro rojo', 'Primera', '300', '28000', '28000', '28000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1867', '15'),
    @(998, '44754', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(999, '44390', 'Zafiro verde', 'Primera', '300', '14000', '14000', '14000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '933', '15'),
    @(1000, '44490', 'Cuatro cascos rojo', 'Primera', '200', '43000', '43000', '43000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '2867', '15'),
    @(1001, '44790', 'Zafiro rojo', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(1002, '44790', 'Zafiro verde', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(1003, '44769', 'Zafiro rojo', 'Primera', '300', '25000', '25000', '25000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1667', '15'),
    @(1004, '44769', 'Zafiro verde', 'Primera', '300', '18000', '18000', '18000', '$/caja 15 kilos', 'Región de Arica y Parinacota', '1200', '15'),
    @(1005, '44890', 'Cuatro cascos verde', 'Primera', '300', '13000', '13000', '13000', '$/caja 18 kilos', 'Región del Maule', '722', '18'),
    @(1006, '44964', 'Cuatro cascos rojo', 'Primera', '150', '12000', '12000', '12000', '$/caja 18 kilos', 'Región del Maule', '667', '18'),
    @(1007, '44964', 'Cuatro cascos verde', 'Primera', '300', '6000', '6000', '6000', '$/caja 18 kilos', 'Región del Maule', '333', '18')
)

foreach ($r in $rowsData) {
    $rowNum = $r[0]
    $ws.Range("D$rowNum").Value = [double]$r[1]
    $ws.Range("H$rowNum").Value = $r[2]
    $ws.Range("I$rowNum").Value = $r[3]
    $ws.Range("J$rowNum").Value = [double]$r[4]
    $ws.Range("K$rowNum").Value = [double]$r[5]
    $ws.Range("L$rowNum").Value = [double]$r[6]
    $ws.Range("M$rowNum").Value = [double]$r[7]
    $ws.Range("N$rowNum").Value = $r[8]
    $ws.Range("O$rowNum").Value = $r[9]
    $ws.Range("P$rowNum").Value = [double]$r[10]
    $ws.Range("Q$rowNum").Value = [double]$r[11]
}

# The two new trailing rows need the constant columns filled in as well,
# since they did not exist in the sheet before.
$ws.Range("A1006").Value = 5
$ws.Range("B1006").Value = "Macroferia Regional de Talca"
$ws.Range("C1006").Value = "Maule"
$ws.Range("E1006").Value = 7
$ws.Range("F1006").Value = 100112002
$ws.Range("G1006").Value = "Pimiento"
$ws.Range("R1006").Value = "Hortaliza"

$ws.Range("A1007").Value = 5
$ws.Range("B1007").Value = "Macroferia Regional de Talca"
$ws.Range("C1007").Value = "Maule"
$ws.Range("E1007").Value = 7
$ws.Range("F1007").Value = 100112002
$ws.Range("G1007").Value = "Pimiento"
$ws.Range("R1007").Value = "Hortaliza"

# Ensure the date-serial columns keep the workbook's date/time number
# format (style index 2 in styles.xml) rather than falling back to
# General for any brand-new cells.
$ws.Range("D922:D1007").NumberFormat = "YYYY-MM-DD HH:MM:SS"
